{"js": "// Spelling/grammar fix in the \"Conclusion\" paragraph of the Glouglou poster:\n//   \"... solide \u00e0 tout personne ...\"        -> \"... solide \u00e0 toute personne ...\"\n//   \"... apprendre ou apponfondire ses ...\"  -> \"... apprendre ou approfondir ses ...\"\nconst body = context.document.body;\n\n// 1) \"tout personne\" -> \"toute personne\" (missing feminine agreement).\nconst toutPersonne = body.search(\"tout personne\", { matchCase: true, matchWholeWord: false });\ntoutPersonne.load(\"items\");\nawait context.sync();\n\nif (toutPersonne.items.length > 0) {\n  toutPersonne.items[0].insertText(\"toute personne\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"apponfondire\" -> \"approfondir\" (misspelled word).\nconst misspelled = body.search(\"apponfondire\", { matchCase: true, matchWholeWord: false });\nmisspelled.load(\"items\");\nawait context.sync();\n\nif (misspelled.items.length > 0) {\n  misspelled.items[0].insertText(\"approfondir\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Spelling/grammar fix in the \"Conclusion\" paragraph of the Glouglou poster:\n#   \"... solide \u00e0 tout personne ...\"        -> \"... solide \u00e0 toute personne ...\"\n#   \"... apprendre ou apponfondire ses ...\"  -> \"... apprendre ou approfondir ses ...\"\n$d = $word.ActiveDocument\n\n# 1) \"tout personne\" -> \"toute personne\" (missing feminine agreement).\n$find1 = $d.Content.Find\n$find1.Text = \"tout personne\"\n$find1.Replacement.Text = \"toute personne\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 2) \"apponfondire\" -> \"approfondir\" (misspelled word).\n$find2 = $d.Content.Find\n$find2.Text = \"apponfondire\"\n$find2.Replacement.Text = \"approfondir\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
